# Fixed the minus in third derivate relation
#
# Column F ("Bisector slope") was computed with the wrong sign/scale on the
# third-derivative term: every value needs to be rescaled by B^4 (B = the
# "Observed wavelength" column) to correct the error.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 479

for ($r = 2; $r -le $lastRow; $r++) {
    $observedWavelength = $ws.Cells.Item($r, 2).Value()   # column B
    $bisectorSlope = $ws.Cells.Item($r, 6).Value()        # column F
    $ws.Cells.Item($r, 6).Value = $bisectorSlope * [Math]::Pow($observedWavelength, 4)
}
